$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 5 (pushes Heatherton .. West Melbourne down by 2)
$ws.Rows.Item(5).Resize(2).Insert()

# Populate the two newly inserted rows with the Glen Waverley exposure sites
$ws.Range("A5").Value = "Glen Waverley"
$ws.Range("B5").Value = "Commonwealth Bank, 28-32 Kingsway, Glen Waverley"
$ws.Range("C5").Value = "1:30pm-2:45pm 9/2/2021"
$ws.Range("D5").Value = "Case attended venue"

$ws.Range("A6").Value = "Glen Waverley"
$ws.Range("B6").Value = "HSBC Bank, 38 Kingsway, Glen Waverley"
$ws.Range("C6").Value = "2:15pm-3:30pm 9/2/2021"
$ws.Range("D6").Value = "Case attended venue"
